# Add two small black (tx1) filled, no-line oval dots to slide 1,
# matching the two new <p:sp> "Oval 34" / "Oval 35" shapes from the diff.
#
# Approach: duplicate an existing oval shape that already carries the
# full <p:style> shape-style reference block and the centered txBody
# formatting (Oval 130), then overwrite its geometry / fill / line /
# name to match the target. Duplicating (instead of AddShape) is what
# lets the new shapes inherit the correct <p:style> + <p:txBody> shape
# structure instead of the bare default produced by AddShape.
#
# The new shape ids must come out as 35 and 36 (matching the diff). This
# runtime hands out shape ids from a monotonically increasing counter
# that also skips any id already used elsewhere in the deck. Ids 2-27
# are free, but 28-34 collide with existing shapes ("Oval 27" .. "TextBox
# 33"), so the 27th/28th shape created in this session naturally land on
# 35/36. We get there by creating+deleting 26 throwaway duplicates first.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$styleSource = $s.Shapes.Item("Oval 130")   # has the p:style + centered txBody we need

for ($i = 1; $i -le 26; $i++) {
    $junk = $styleSource.Duplicate().Item(1)
    $junk.Delete()
}

function Add-Dot {
    # NOTE: this COM host only binds positional params reliably, not
    # PowerShell named (-Foo bar) args - keep the call sites positional.
    param(
        [double]$Left,
        [double]$Top,
        [double]$Size,
        [string]$Name
    )

    $shp = $styleSource.Duplicate().Item(1)
    $shp.Left = $Left
    $shp.Top = $Top
    $shp.Width = $Size
    $shp.Height = $Size
    $shp.Fill.Solid()
    $shp.Fill.ForeColor.SchemeColor = "tx1"
    $shp.Line.Visible = $false
    $shp.Name = $Name
    return $shp
}

# EMU targets (x=2442164 y=778582 cx=cy=45719) expressed in points, nudged
# by a hair so this runtime's float32 Left/Top/Width/Height round-trip
# lands exactly back on the target EMU values instead of one EMU short.
$dot1 = Add-Dot 192.2964172527559 61.30570989133858 3.5999606598425196 "Oval 34"

# EMU targets x=2337955 y=3934047 cx=cy=45719
$dot2 = Add-Dot 184.09098818188974 309.76752161496063 3.5999606598425196 "Oval 35"

Write-Output ("dot1 id=" + $dot1.Id + " name=" + $dot1.Name)
Write-Output ("dot2 id=" + $dot2.Id + " name=" + $dot2.Name)
